$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Vsilva" to "Vsilva2"
$ws.Name = "Vsilva2"

# Move the selection/active cell to M17
$ws.Range("M17").Select()
